$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "'C21, C22, CH1_C23, CH1_C24, CH2_C23, CH2_C24, CH3_C23, CH3_C24, CH4_C23, CH4_C24, CH5_C23, CH5_C24"
$ws.Range("A5").Value = "'CH1_J2, CH2_J2, CH3_J2, CH4_J2, CH5_J2"
$ws.Range("A6").Value = "'CH1_R7, CH1_R8, CH2_R7, CH2_R8, CH3_R7, CH3_R8, CH4_R7, CH4_R8, CH5_R7, CH5_R8, R3, R4, R5, R6"
$ws.Range("A7").Value = "'CH1_U9, CH2_U9, CH3_U9, CH4_U9, CH5_U9, U7"
